$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet (and update title text) from 2022-08-02 to 2022-08-03
$ws.Name = "Through 2022-08-03"

# Update the August label cell
$ws.Range("A9").Value = "August (through 08-03)"

# Correction to 2022 value in row 7 (June)
$ws.Range("I7").Value = 142

# Update August row data (row 9)
$ws.Range("B9").Value = 2
$ws.Range("D9").Value = 7
$ws.Range("E9").Value = 7
$ws.Range("G9").Value = 20
$ws.Range("H9").Value = 14
$ws.Range("I9").Value = 16

# Update Total row data (row 10)
$ws.Range("B10").Value = 164
$ws.Range("D10").Value = 472
$ws.Range("E10").Value = 432
$ws.Range("G10").Value = 641
$ws.Range("H10").Value = 924
$ws.Range("I10").Value = 986
